$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the percent number format that was mistakenly applied to the ticker
#        column (A2:A5) -- remove it, leaving the border / bold font / center+top
#        alignment untouched.
$ws.Range("A2:A5").NumberFormat = "general"

# --- 2. Update the ticker symbols (row labels)
$ws.Range("A2").Value2 = "UAL"
$ws.Range("A3").Value2 = "AAL"
$ws.Range("A4").Value2 = "LUV"
$ws.Range("A5").Value2 = "SAVE"

# --- 3. Update the column headers (row 1)
$ws.Range("B1").Value2 = "Long-term debt"
$ws.Range("C1").Value2 = "Gross profit"
$ws.Range("D1").Value2 = "Price to free cash flows ratio"
$ws.Range("E1").Value2 = "Price to book ratio"
$ws.Range("F1").Value2 = "Price to sales ratio"

# --- 4. Move the percent-style formatting that used to live on columns B, C and F
#        onto column E, then strip it from B, C and F (those columns have no other
#        formatting of their own, so resetting to the default/general format leaves
#        plain, unstyled cells).
$pctFormat = $ws.Range("F2").NumberFormat
$ws.Range("E2:E5").NumberFormat = $pctFormat
$ws.Range("B2:B5").NumberFormat = "general"
$ws.Range("C2:C5").NumberFormat = "general"
$ws.Range("F2:F5").NumberFormat = "general"

# --- 5. Write the new data values
$ws.Range("B2").Value2 = 28283
$ws.Range("C2").Value2 = 15082
$ws.Range("D2").Value2 = 9.789629593023255
$ws.Range("E2").Value2 = 1.770253495142111
$ws.Range("F2").Value2 = 0.2715530664553442

$ws.Range("B3").Value2 = 32389
$ws.Range("C3").Value2 = 12439
$ws.Range("D3").Value2 = -11.28565948158254
$ws.Range("E3").Value2 = -1.426519813760993
$ws.Range("F3").Value2 = 0.1689242286251046

$ws.Range("B4").Value2 = 8046
$ws.Range("C4").Value2 = 4503
$ws.Range("D4").Value2 = -127.9891666666667
$ws.Range("E4").Value2 = 1.868280153457472
$ws.Range("F4").Value2 = 0.8384273956496179

$ws.Range("B5").Value2 = 3200.376
$ws.Range("C5").Value2 = 799.317
$ws.Range("D5").Value2 = -5.996743227559657
$ws.Range("E5").Value2 = 1.347926148998728
$ws.Range("F5").Value2 = 0.4179721086163079

# --- 6. Remove the now-unused columns G:J (the table shrank from 10 to 6 columns)
$ws.Range("G1:J1").EntireColumn.Delete() | Out-Null

Write-Output "done"
